$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New "Email" column header in D1 (reuses existing shared string "Email")
$ws.Range("D1").Value = "Email"

# Row 15: Wordpress Blog (main account)
# Write in A, C, B, D order to reproduce the original shared-string insertion order
$ws.Range("A15").Value = "Wordpress Blog"
$ws.Range("C15").Value = "Peoplespaceoc2"
$ws.Range("B15").Value = "womencoders-admin"
$ws.Range("D15").Value = "womencoders@gmail.com"

# Row 16: Wordpress Blog - Angela
$ws.Range("A16").Value = "Wordpress Blog - Angela"
$ws.Range("B16").Value = "angelal4"
$ws.Range("D16").Value = "angelgirl2272@gmail.com"

# Row 17: Wordpress Blog - Laurie
$ws.Range("A17").Value = "Wordpress Blog - Laurie"
$ws.Range("B17").Value = "laurie415"
$ws.Range("D17").Value = "ldxtran@gmail.com"

# Hyperlinks on the new Email cells (mailto links), styled like the rest of the sheet
$ws.Hyperlinks.Add($ws.Range("D15"), "mailto:womencoders@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D16"), "mailto:angelgirl2272@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D17"), "mailto:ldxtran@gmail.com")

$ws.Range("D15").Style = "Hyperlink"
$ws.Range("D16").Style = "Hyperlink"
$ws.Range("D17").Style = "Hyperlink"

# Restore the selection as left by the author
$ws.Range("D12").Select()
